$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.176.45'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -4.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.487.41'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.66%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.40'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.33'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -6.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.33%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.574'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.517.33'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0999'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -4.59%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.54'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.926.84'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.03'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -6.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.072.69'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.62%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.502.64'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -3.06%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -5.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.32'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.76'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.36'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.442'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -10.56%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.607.51'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -3.40%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.83'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -4.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.90'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -6.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0779'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -7.33%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.86%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.50%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.30%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.16'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.45'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +1.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.54'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.45'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -9.00%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -10.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.93'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '36.78'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.68'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -6.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.821'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -9.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.994'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.596'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.78'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.46'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0932'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.61'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -5.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0518'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -5.45%  '
